$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6.224379325864566
$ws.Range("G2").Value = 8.928571428571429
$ws.Range("F3").Value = 6.224379325864566
$ws.Range("D6").Value = "Master-Mañanas"
$ws.Range("G10").Value = 10.71428571428572
$ws.Range("G13").Value = 80.35714285714286
